$d = $word.ActiveDocument

# Start from the very end of the document and append:
#   1) a blank separator paragraph
#   2) a bold "Julieta Simos" paragraph
#   3) her comment paragraph (regular formatting)
# mirroring the pattern already used for every other contributor in the file.

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

$sepPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$sepPara.Range.InsertParagraphAfter()

$namePara = $d.Paragraphs.Item($d.Paragraphs.Count)
$namePara.Range.InsertParagraphAfter()

$namePara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$namePara.Range.Text = "Julieta Simos"
$namePara.Range.Font.Bold = 1

$commentPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$commentPara.Range.Text = "La realizacion de este trabajo practico en grupo fue interesante, una experiencia necesaria para poder entender e incorporar los conocimientos vistos, pero ademas de esto fue gratificante,  al poder trabajar en equipo el trabajo se hace mas llevadero y es mas dificil que se pierda la motivacion. El utilizar nuevas herramientas para organizar mejor las tareas nos da una idea de como podria ser en un futuro la organizacion y boceto de nuestro propio proyecto. Como todo nuevo conocimiento tiene su parte dificil, el conocer y utilizar los comandos de git fue una tarea a veces complicada, el miedo a cometer un error y no poder resolverlo, a pesar de los errores o conflictos sucedidos se pudieron resolver buscando e investigando sobre el tema, conversando y ayudandonos entre si, con lo que finalmente pudimos terminar el trabajo practico."
